$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1485
    $ws.Range("F6").Value = 33
    $ws.Range("F9").Value = 269
}
